# Presentation for the experiment added
# - Convert the Sheet1 chart from a line chart to a clustered column (bar) chart
# - Re-color/size the chart to the new bar-chart defaults (fill instead of stroke, gap width, legend on the right)
# - Shrink/reposition the chart's anchor on the worksheet
# - Update the worksheet's active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart

# --- Chart type: line -> clustered column ---
$chart.ChartType = 51   # xlColumnClustered (barDir="col", grouping="clustered")

# --- Per-series tweaks: bars are now filled with the accent color instead of
#     being drawn as a colored stroke, and bars don't invert color on negative
#     values. Also flatten the leftover stroke so it doesn't show as an outline. ---
$ser1 = $chart.SeriesCollection(1)
$ser2 = $chart.SeriesCollection(2)

$ser1.Interior.Color = 0xBD814F   # accent1 (4F81BD) - COM ColorFormat is BGR-ordered
$ser2.Interior.Color = 0x4D50C0   # accent2 (C0504D) - COM ColorFormat is BGR-ordered

$ser1.Border.Weight = 0
$ser2.Border.Weight = 0

$ser1.InvertIfNegative = $false
$ser2.InvertIfNegative = $false

# --- Bar chart group spacing ---
$chart.ChartGroups(1).GapWidth = 150

# --- Legend moves from bottom to the right ---
$chart.Legend.Position = -4152   # xlLegendPositionRight

# --- Resize/reposition the chart object on the sheet (was ~17 cols x 32 rows,
#     now smaller at ~14.3 cols x 29.5 rows) while keeping the top-left corner. ---
$chartObj.Width = 787.625
$chartObj.Height = 433.5

# --- Update the worksheet's active selection ---
$null = $ws.Range("Y23").Select()
